$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("I3").Value = 1.33
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 13
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53

# Row 6 updates
$ws.Range("G6").Value = 1.52
$ws.Range("H6").Value = 3.5
$ws.Range("I6").Value = 6.9
$ws.Range("J6").Value = 2.02
$ws.Range("K6").Value = 2.12
$ws.Range("L6").Value = 6.5
$ws.Range("N6").Value = 6.65
$ws.Range("O6").Value = 1.34
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.65
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.52
$ws.Range("U6").Value = 2.05
$ws.Range("V6").Value = 1.62
$ws.Range("W6").Value = 5.4
$ws.Range("X6").Value = 6.2
$ws.Range("Y6").Value = 8.25
$ws.Range("Z6").Value = 10.25
$ws.Range("AA6").Value = 13.5
$ws.Range("AC6").Value = 8
$ws.Range("AD6").Value = 7.1
$ws.Range("AE6").Value = 20
$ws.Range("AF6").Value = 110
$ws.Range("AH6").Value = 15.5
$ws.Range("AI6").Value = 45
$ws.Range("AJ6").Value = 22
$ws.Range("AK6").Value = 175
$ws.Range("AM6").Value = 80
$ws.Range("AN6").Value = 3.2
$ws.Range("AO6").Value = 7.1
$ws.Range("AP6").Value = 18
$ws.Range("AQ6").Value = 23
$ws.Range("AR6").Value = 55
$ws.Range("AT6").Value = 2.47
$ws.Range("AU6").Value = 7.9
$ws.Range("AW6").Value = 7.9
$ws.Range("AY6").Value = 40
$ws.Range("BA6").Value = 300

$wb.Save()
